$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 299.19354
$ws.Range("I33").Value = 264.65518
$ws.Range("K33").Value = 264.65518
$ws.Range("M33").Value = -35.65517999999997
$ws.Range("H43").Value = 1826.6666
$ws.Range("I43").Value = 1850
$ws.Range("J43").Value = 1780
$ws.Range("K43").Value = 1850
$ws.Range("L43").Value = 1780
$ws.Range("M43").Value = -1781
$ws.Range("N43").Value = -1918
$ws.Range("H64").Value = 3894.9
$ws.Range("I64").Value = 3712.375
$ws.Range("J64").Value = 4016.5833
$ws.Range("K64").Value = 3712.375
$ws.Range("L64").Value = 4016.5833
$ws.Range("M64").Value = -3464.375
$ws.Range("N64").Value = -4512.5833
$ws.Range("H67").Value = 3894.9
$ws.Range("I67").Value = 3712.375
$ws.Range("J67").Value = 4016.5833
$ws.Range("K67").Value = 3712.375
$ws.Range("L67").Value = 4016.5833
$ws.Range("M67").Value = -2854.375
$ws.Range("N67").Value = -5732.5833
$ws.Range("H74").Value = 2500
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1564
$ws.Range("N74").ClearContents() | Out-Null
$ws.Range("H77").Value = 2500
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -7820
$ws.Range("N77").ClearContents() | Out-Null
$ws.Range("H107").Value = 678.94116
$ws.Range("I107").Value = 604.3
$ws.Range("J107").Value = 785.5714
$ws.Range("K107").Value = 604.3
$ws.Range("L107").Value = 785.5714
$ws.Range("M107").Value = 1315.7
$ws.Range("N107").Value = -4625.5714
$ws.Range("H129").Value = 843.28125
$ws.Range("J129").Value = 848.68854
$ws.Range("L129").Value = 2546.06562
$ws.Range("N129").Value = -12546.06562
$ws.Range("H132").Value = 2662.2334
$ws.Range("I132").Value = 2880.0386
$ws.Range("J132").Value = 1246.5
$ws.Range("K132").Value = 8640.1158
$ws.Range("L132").Value = 3739.5
$ws.Range("M132").Value = -6110.1158
$ws.Range("N132").Value = -8799.5
$ws.Range("H138").Value = 2540.5806
$ws.Range("I138").Value = 1603.7894
$ws.Range("J138").Value = 4023.8333
$ws.Range("K138").Value = 4811.3682
$ws.Range("L138").Value = 12071.4999
$ws.Range("M138").Value = 328.6318000000001
$ws.Range("N138").Value = -22351.4999
$ws.Range("H141").Value = 3398.125
$ws.Range("J141").Value = 4450
$ws.Range("L141").Value = 13350
$ws.Range("N141").Value = -23710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 6250
$ws.Range("I31").Value = 6250
$ws.Range("K31").Value = 6250
$ws.Range("M31").Value = -5956
$ws.Range("H32").Value = 1560.1765
$ws.Range("I32").Value = 1333.4833
$ws.Range("J32").Value = 3260.375
$ws.Range("K32").Value = 1333.4833
$ws.Range("L32").Value = 3260.375
$ws.Range("M32").Value = -1046.4833
$ws.Range("N32").Value = -3834.375
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2533
$ws.Range("H93").Value = 29400
$ws.Range("J93").Value = 29400
$ws.Range("L93").Value = 29400
$ws.Range("N93").Value = -34392

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1666.5555
$ws.Range("I99").Value = 1373.091
$ws.Range("J99").Value = 2957.8
$ws.Range("K99").Value = 1373.091
$ws.Range("L99").Value = 2957.8
$ws.Range("M99").Value = 124.9090000000001
$ws.Range("N99").Value = -5953.8
$ws.Range("H134").Value = 4263.952
$ws.Range("I134").Value = 4263.952
$ws.Range("K134").Value = 12791.856
$ws.Range("M134").Value = -10256.856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1188.4445
$ws.Range("I16").Value = 1170.8572
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 1170.8572
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -883.8571999999999
$ws.Range("N16").Value = -1824
$ws.Range("H99").Value = 4355.304
$ws.Range("I99").Value = 3338.9412
$ws.Range("J99").Value = 7235
$ws.Range("K99").Value = 3338.9412
$ws.Range("L99").Value = 7235
$ws.Range("M99").Value = -1840.9412
$ws.Range("N99").Value = -10231
$ws.Range("H113").Value = 1188.4445
$ws.Range("I113").Value = 1170.8572
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 1170.8572
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = 999.1428000000001
$ws.Range("N113").Value = -5590
$ws.Range("H126").Value = 4355.304
$ws.Range("I126").Value = 3338.9412
$ws.Range("J126").Value = 7235
$ws.Range("K126").Value = 10016.8236
$ws.Range("L126").Value = 21705
$ws.Range("M126").Value = -7546.8236
$ws.Range("N126").Value = -26645
$ws.Range("H134").Value = 1097.1177
$ws.Range("I134").Value = 903.7143
$ws.Range("J134").Value = 1999.6666
$ws.Range("K134").Value = 2711.1429
$ws.Range("L134").Value = 5998.9998
$ws.Range("M134").Value = -176.1428999999998
$ws.Range("N134").Value = -11068.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 5447.5557
$ws.Range("J96").Value = 5447.5557
$ws.Range("L96").Value = 16342.6671
$ws.Range("N96").Value = -20460.6671
$ws.Range("H131").Value = 150052.31
$ws.Range("I131").Value = 732
$ws.Range("J131").Value = 162094.28
$ws.Range("K131").Value = 2196
$ws.Range("L131").Value = 486282.84
$ws.Range("M131").Value = 2844
$ws.Range("N131").Value = -496362.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4610.8887
$ws.Range("I70").Value = 4428.4287
$ws.Range("J70").Value = 4727
$ws.Range("K70").Value = 4428.4287
$ws.Range("L70").Value = 4727
$ws.Range("M70").Value = -4158.4287
$ws.Range("N70").Value = -5267
$ws.Range("H73").Value = 4610.8887
$ws.Range("I73").Value = 4428.4287
$ws.Range("J73").Value = 4727
$ws.Range("K73").Value = 4428.4287
$ws.Range("L73").Value = 4727
$ws.Range("M73").Value = -3492.4287
$ws.Range("N73").Value = -6599
$ws.Range("H113").Value = 3045.1
$ws.Range("I113").Value = 2207.2856
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2207.2856
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -37.28560000000016
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5805.737
$ws.Range("I61").Value = 3250.5
$ws.Range("J61").Value = 7664.091
$ws.Range("K61").Value = 3250.5
$ws.Range("L61").Value = 7664.091
$ws.Range("M61").Value = -3048.5
$ws.Range("N61").Value = -8068.091
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents() | Out-Null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents() | Out-Null
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents() | Out-Null
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents() | Out-Null
$ws.Range("H113").Value = 5805.737
$ws.Range("I113").Value = 3250.5
$ws.Range("J113").Value = 7664.091
$ws.Range("K113").Value = 3250.5
$ws.Range("L113").Value = 7664.091
$ws.Range("M113").Value = -1080.5
$ws.Range("N113").Value = -12004.091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 3038
$ws.Range("J18").Value = 3038
$ws.Range("L18").Value = 3038
$ws.Range("N18").Value = -3384
$ws.Range("H19").Value = 2000
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -2348
$ws.Range("H81").Value = 1983.3334
$ws.Range("I81").Value = 1983.3334
$ws.Range("K81").Value = 3966.6668
$ws.Range("M81").Value = -2905.6668
$ws.Range("H84").Value = 1983.3334
$ws.Range("I84").Value = 1983.3334
$ws.Range("K84").Value = 19833.334
$ws.Range("M84").Value = -14529.334
